$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 10, 11, 13, 14)

foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-19 08:21:55"

    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-19 08:21:50"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-19 08:21:55"
}
